$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("D2893_2_bg_detlim")
$ws2 = $wb.Worksheets.Item("D2893_3_bg_apf_detlim")

# Row data layout: columns B..Z (2..26) hold values for rows 2, 5 and 8 on each sheet.

$sheet1Row2 = @(0.015,0.017,0.017,0.021,0.014,0.03,0.014,0.018,0.016,0.015,0.015,0.015,0.015,0.013,0.014,0.015,0.029,0.033,0.026,0.027,0.029,0.019,0.007,0.013,0.033)
$sheet1Row5 = @(0.015,0.017,0.017,0.021,0.014,0.03,0.014,0.018,0.016,0.015,0.015,0.015,0.015,0.013,0.014,0.015,0.029,0.033,0.026,0.027,0.029,0.019,0.007,0.013,0.033)
$sheet1Row8 = @(0.022,0.025,0.024,0.031,0.02,0.043,0.02,0.026,0.024,0.022,0.021,0.021,0.022,0.018,0.02,0.021,0.041,0.048,0.038,0.039,0.042,0.028,0.01,0.018,0.048)

$sheet2Row2 = @(0.018,0.021,0.02,0.026,0.016,0.036,0.017,0.022,0.02,0.018,0.018,0.018,0.018,0.015,0.017,0.018,0.034,0.04,0.032,0.033,0.035,0.023,0.008,0.015,0.04)
$sheet2Row5 = @(0.018,0.021,0.02,0.026,0.016,0.036,0.017,0.022,0.02,0.018,0.018,0.018,0.018,0.015,0.017,0.018,0.034,0.04,0.032,0.033,0.035,0.023,0.008,0.015,0.04)
$sheet2Row8 = @(0.026,0.03,0.029,0.037,0.024,0.052,0.024,0.032,0.028,0.026,0.025,0.026,0.026,0.022,0.024,0.025,0.05,0.057,0.046,0.047,0.051,0.034,0.011,0.022,0.057)

function Set-RowValues($ws, $rowNum, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B is index 2
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}

Set-RowValues $ws1 2 $sheet1Row2
Set-RowValues $ws1 5 $sheet1Row5
Set-RowValues $ws1 8 $sheet1Row8

Set-RowValues $ws2 2 $sheet2Row2
Set-RowValues $ws2 5 $sheet2Row5
Set-RowValues $ws2 8 $sheet2Row8
